$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new "Save" header in H1, copying the formatting (bold, border,
# centered alignment) used by the other header cells (e.g. G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new "Save" column values (H2:H6), all zero, matching the
# plain (unstyled) numeric formatting used by the other data columns.
$ws.Range("H2:H6").Value = 0
